$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the unit prices in column D
$ws.Range("D31").Value = 3789.125
$ws.Range("D32").Value = 4996.397
$ws.Range("D33").Value = 6321.77
$ws.Range("D34").Value = 7128.801
